$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "37.055.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.82%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.091.60"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.66%  "

$ws.Range("E4").Value = "  -0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "249.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.650"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.07%  "

$ws.Range("E7").Value = "  -0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.90"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +8.99%  "

$ws.Range("E9").Value = "  +4.55%  "

$ws.Range("E10").Value = "  -3.36%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.50%  "

$ws.Range("E12").Value = "  +4.78%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.47%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.358.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.08%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.825"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.50%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.091.50"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.78%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.80%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "36.958.23"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -4.41%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0820"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.79%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.29"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "239.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.71%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.19"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.05%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.09%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "169.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.58%  "

$ws.Range("E27").Value = "  +2.07%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.54"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -5.43%  "

$ws.Range("E30").Value = "  -5.90%  "

$ws.Range("E31").Value = "  +16.46%  "

$ws.Range("E32").Value = "  -2.22%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0604"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "20.77"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -8.90%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0909"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.38%  "

$ws.Range("E36").Value = "  -0.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.30"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.07"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.24%  "

$ws.Range("E39").Value = "  -3.86%  "

$ws.Range("E40").Value = "  -9.11%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "17.57"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.96%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0223"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.77%  "

$ws.Range("E43").Value = "  +1.31%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.89"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.24%  "

$ws.Range("E45").Value = "  -5.15%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0883"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.91%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.99"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.308.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.63%  "

$ws.Range("E49").Value = "  +4.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.277.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.57%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -5.72%  "
